# Updates the "cryptos" price/volume table in place.
#
# Columns D (Price) and E (Volume(1h)) are stored as plain text in the
# source workbook, even though many of the Price values look like plain
# decimal numbers (e.g. "61.10"). Assigning such a string straight to
# Range.Value lets Excel's usual literal-entry parsing turn it into a
# real number (dropping the trailing zero, e.g. 61.10 -> 61.1), which
# would not match the source data. To keep those cells as text we enter
# them with a leading apostrophe (the standard "force text" prefix) and
# then copy the row's plain (unstyled) cell style back onto the cell so
# no left-over "quote prefix" formatting remains on it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.942.24"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").Value = "2.650.85"
$ws.Range("E3").Value = "  +3.67%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'514.98"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +1.86%  "

$ws.Range("D6").Value = "'144.31"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  -0.46%  "

$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("D8").Value = "'0.568"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "  +2.06%  "

$ws.Range("D9").Value = "2.680.69"
$ws.Range("E9").Value = "  +5.13%  "

$ws.Range("D10").Value = "'6.32"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  +2.03%  "

$ws.Range("D11").Value = "'0.106"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  +3.92%  "

$ws.Range("D12").Value = "'0.337"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "  +1.43%  "

$ws.Range("E13").Value = "  -1.20%  "

$ws.Range("D14").Value = "3.115.56"
$ws.Range("E14").Value = "  +3.49%  "

$ws.Range("D15").Value = "58.933.06"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").Value = "'21.12"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = "  +2.83%  "

$ws.Range("E17").Value = "  +2.26%  "

$ws.Range("D18").Value = "2.674.26"
$ws.Range("E18").Value = "  +4.47%  "

$ws.Range("D19").Value = "'4.55"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("D20").Value = "'341.74"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  +2.25%  "

$ws.Range("D21").Value = "'10.46"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "  +3.85%  "

$ws.Range("D22").Value = "'6.14"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  +3.47%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").Value = "'61.10"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  +2.64%  "

$ws.Range("D25").Value = "'0.421"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "  +3.28%  "

$ws.Range("D26").Value = "2.761.57"
$ws.Range("E26").Value = "  +3.10%  "

$ws.Range("D27").Value = "'0.993"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("E28").Value = "  +3.75%  "

$ws.Range("D29").Value = "0.0₃0812"
$ws.Range("E29").Value = "  +4.49%  "

$ws.Range("D30").Value = "'7.15"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  +4.32%  "

$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("E32").Value = "  +9.05%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'18.95"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = "  +2.00%  "

$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "'1.58"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  +2.76%  "

$ws.Range("D35").Value = "'149.20"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  +0.36%  "

$ws.Range("E36").Value = "  +13.35%  "

$ws.Range("E37").Value = "  +4.41%  "

$ws.Range("E38").Value = "  +3.95%  "

$ws.Range("D39").Value = "'0.856"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = "  +4.85%  "

$ws.Range("D40").Value = "'36.64"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  +2.18%  "

$ws.Range("D41").Value = "'3.69"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  +4.60%  "

$ws.Range("E42").Value = "  +1.96%  "

$ws.Range("D43").Value = "'283.38"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  +1.00%  "

$ws.Range("D44").Value = "'0.620"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  +2.47%  "

$ws.Range("E45").Value = "  -0.42%  "

$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").Value = "'19.54"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  +4.83%  "

$ws.Range("D48").Value = "'0.0535"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("E49").Value = "  +1.48%  "

$ws.Range("D50").Value = "'4.72"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  +4.84%  "

$ws.Range("D51").Value = "'10.28"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  -0.50%  "
